$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metro_budget")

# ---------------------------------------------------------------------
# Question 9 (rows 98:100) - fill in the Department / Pct lookup table,
# mirroring the "Question 8" block above (rows 91:93) but using
# INDEX/MATCH instead of XLOOKUP, per the finished workbook.
# ---------------------------------------------------------------------

# Row 98 - anchor formulas for the shared groups below
$ws.Range("B98").Formula = '=INDEX($A$1:$A$52,MATCH($B$96,INDEX($A$1:$P$52,,MATCH(_xlfn.CONCAT($A98,"_rank"),$A$1:$P$1,0)),0))'
$ws.Range("C98").Formula = '=INDEX($A$1:$P$52,MATCH($B98, $A$1:$A$52,0),MATCH(_xlfn.CONCAT($A98,"_diff_pct"),$A$1:$P$1,0))'
$ws.Range("D98").Formula = '=INDEX($A$1:$A$52,MATCH($D$96,INDEX($A$1:$P$52,,MATCH(_xlfn.CONCAT($A98,"_rank"),$A$1:$P$1,0)),0))'
$ws.Range("E98").Formula = '=INDEX($A$1:$P$52,MATCH($D98, $A$1:$A$52,0),MATCH(_xlfn.CONCAT($A98,"_diff_pct"),$A$1:$P$1,0))'
$ws.Range("F98").Formula = '=INDEX($A$1:$A$52,MATCH($F$96,INDEX($A$1:$P$52,,MATCH(_xlfn.CONCAT($A98,"_rank"),$A$1:$P$1,0)),0))'
$ws.Range("G98").Formula = '=INDEX($A$1:$P$52,MATCH($F98, $A$1:$A$52,0),MATCH(_xlfn.CONCAT($A98,"_diff_pct"),$A$1:$P$1,0))'

# Rows 99:100 - fill down (creates the shared formula groups)
$ws.Range("B99:B100").Formula = '=INDEX($A$1:$A$52,MATCH($B$96,INDEX($A$1:$P$52,,MATCH(_xlfn.CONCAT($A99,"_rank"),$A$1:$P$1,0)),0))'
$ws.Range("C99:C100").Formula = '=INDEX($A$1:$P$52,MATCH($B99, $A$1:$A$52,0),MATCH(_xlfn.CONCAT($A99,"_diff_pct"),$A$1:$P$1,0))'
$ws.Range("D99:D100").Formula = '=INDEX($A$1:$A$52,MATCH($D$96,INDEX($A$1:$P$52,,MATCH(_xlfn.CONCAT($A99,"_rank"),$A$1:$P$1,0)),0))'
$ws.Range("E99:E100").Formula = '=INDEX($A$1:$P$52,MATCH($D99, $A$1:$A$52,0),MATCH(_xlfn.CONCAT($A99,"_diff_pct"),$A$1:$P$1,0))'
$ws.Range("F99:F100").Formula = '=INDEX($A$1:$A$52,MATCH($F$96,INDEX($A$1:$P$52,,MATCH(_xlfn.CONCAT($A99,"_rank"),$A$1:$P$1,0)),0))'
$ws.Range("G99:G100").Formula = '=INDEX($A$1:$P$52,MATCH($F99, $A$1:$A$52,0),MATCH(_xlfn.CONCAT($A99,"_diff_pct"),$A$1:$P$1,0))'

# The diff/pct columns (E, G) pick up the 0.00% number format (column C
# already carries it from the original template).
$ws.Range("E98:E100").NumberFormat = "0.00%"
$ws.Range("G98:G100").NumberFormat = "0.00%"

# ---------------------------------------------------------------------
# Reposition the embedded chart: nudged up-and-right slightly (same size).
# ---------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$co.Top = $co.Top() - 9.75
$co.Left = $co.Left() + 3.75

# ---------------------------------------------------------------------
# Selection / scroll state left by the author when they saved.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 80
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F104").Select()
